# This edit re-orders the 28 trial rows (rows 2-29) of the input-list sheet
# into the final randomized presentation order used for the experiment run,
# renumbers trial_total (col F) to start at 297 instead of 327, and swaps the
# placeholder catch-trial stimulus for the real one now that the list is final.
#
# Strategy: snapshot every data row's movable columns (2-29, columns G-S) as
# they exist today, then write each destination row from the snapshot of its
# source row, following the de-duplicated "after <- before" row mapping below.
# Columns A-E (subject_id..trial_block) are untouched - trial_block is just
# the 1..28 row position. Column F (trial_total) is renumbered separately.
# The old catch row (originally row 9, stimuli/catch_18.jpg) is retired; a
# brand new catch row (stimuli/catch_01.jpg) lands at row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 29
$lastCol = 19   # column S

# --- 1. snapshot current rows (columns G..S only - "A".."E" never move,
#        and "F" is recomputed from scratch in step 4) ------------------------
$firstMovedCol = 7   # column G
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstMovedCol; $c -le $lastCol; $c++) {
        $rowVals += , $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# --- 2. destination row -> source row mapping (data movement only) -----------
# Row 18 has no source: it becomes the new catch trial.
$rowMap = @{
    2  = 25
    3  = 4
    4  = 27
    5  = 15
    6  = 10
    7  = 18
    8  = 21
    9  = 3
    10 = 19
    11 = 29
    12 = 7
    13 = 26
    14 = 17
    15 = 16
    16 = 6
    17 = 14
    18 = 0   # new catch row, no source
    19 = 22
    20 = 2
    21 = 13
    22 = 24
    23 = 20
    24 = 23
    25 = 8
    26 = 11
    27 = 5
    28 = 28
    29 = 12
}

# --- 3. write every destination row from its source snapshot (cols G..S) -----
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $rowMap[$r]
    if ($src -ne 0) {
        $vals = $snapshot[$src]
        for ($c = $firstMovedCol; $c -le $lastCol; $c++) {
            $ws.Cells.Item($r, $c).Value = $vals[$c - $firstMovedCol]
        }
    }
}

# --- 4. renumber trial_total (col F, 297..324) --------------------------------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = 297 + ($r - $firstRow)
}

# --- 5. rebuild row 18 as the new catch trial ---------------------------------
# category / cond_cat / conceptual.. columns are blank on catch rows.
$ws.Cells.Item(18, 7).Value = "living_rooms"      # G: target_cat
$ws.Cells.Item(18, 8).ClearContents()              # H: category
$ws.Cells.Item(18, 9).ClearContents()              # I: cond_cat
$ws.Cells.Item(18, 10).Value = "catch"             # J: cond_mem
$ws.Cells.Item(18, 11).Value = "f"                 # K: correct_answer
$ws.Cells.Item(18, 12).Value = "stimuli/catch_01.jpg"  # L: stimulus
$ws.Cells.Item(18, 13).ClearContents()             # M: conceptual
$ws.Cells.Item(18, 14).ClearContents()             # N: perceptual
$ws.Cells.Item(18, 15).ClearContents()             # O: typicality
$ws.Cells.Item(18, 16).ClearContents()             # P: n
$ws.Cells.Item(18, 17).ClearContents()             # Q: p_typicality
$ws.Cells.Item(18, 18).ClearContents()             # R: p_conceptual
$ws.Cells.Item(18, 19).ClearContents()             # S: p_perceptual
